# "load all sheets or tables into tabs" - add a second worksheet ("Another
# Sheet") holding a small Name/title/status table, alongside the existing
# "sampledata" sheet.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Duplicate the existing sheet so the new tab starts out with the workbook's
# normal formatting/drawing plumbing already wired up (font/style reuse,
# an attached drawing part, etc.), then wipe its contents and rename it.
$ws1.Copy([System.Reflection.Missing]::Value, $ws1)
$newSheet = $wb.Worksheets.Item(2)
$newSheet.Name = "Another Sheet"
$newSheet.Cells.Clear()

# The little roster table to load into the new tab.
$data = @(
    @("Name", "title", "status"),
    @("Bob", "Doctor", "Active"),
    @("Mike", "Technician", "Inactive"),
    @("Adam", "Driver", "Active"),
    @("Kelly", "Actor", "Inactive")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    for ($j = 0; $j -lt $data[$i].Count; $j++) {
        $newSheet.Cells.Item($i + 1, $j + 1).Value = $data[$i][$j]
    }
}

# Reuse the original sheet's cell styling (font) instead of minting a new
# style entry in styles.xml.
$ws1.Range("A1:C1").Copy()
$newSheet.Range("A1:C5").PasteSpecial(-4122)  # xlPasteFormats

# Leave the original sheet active/selected, matching the workbook's prior
# (default) view state.
$ws1.Activate() | Out-Null
$ws1.Range("A1").Select() | Out-Null
